$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — update F2 and F4
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1012
$ws1.Range("F4").Value = 437

# Sheet "全部类型" (All Types) — update F2 and F4 (mirrors the exhibition sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1012
$ws4.Range("F4").Value = 437
